$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "34÷4=8, 2"
$t.Cell(1, 2).Range.Text = "18÷8=2, 2"
$t.Cell(1, 3).Range.Text = "92÷4=23, 0"
$t.Cell(1, 4).Range.Text = "71÷4=17, 3"
$t.Cell(1, 5).Range.Text = "41÷3=13, 2"
$t.Cell(5, 1).Range.Text = "88÷2=44, 0"
$t.Cell(5, 2).Range.Text = "25÷2=12, 1"
$t.Cell(5, 3).Range.Text = "27÷4=6, 3"
$t.Cell(5, 4).Range.Text = "35÷4=8, 3"
$t.Cell(5, 5).Range.Text = "25÷8=3, 1"
$t.Cell(9, 1).Range.Text = "24÷7=3, 3"
$t.Cell(9, 2).Range.Text = "81÷9=9, 0"
$t.Cell(9, 3).Range.Text = "88÷5=17, 3"
$t.Cell(9, 4).Range.Text = "25÷8=3, 1"
$t.Cell(9, 5).Range.Text = "19÷9=2, 1"
$t.Cell(13, 1).Range.Text = "75÷5=15, 0"
$t.Cell(13, 2).Range.Text = "98÷3=32, 2"
$t.Cell(13, 3).Range.Text = "90÷4=22, 2"
$t.Cell(13, 4).Range.Text = "78÷5=15, 3"
$t.Cell(13, 5).Range.Text = "50÷3=16, 2"
$t.Cell(17, 1).Range.Text = "71÷7=10, 1"
$t.Cell(17, 2).Range.Text = "55÷3=18, 1"
$t.Cell(17, 3).Range.Text = "38÷6=6, 2"
$t.Cell(17, 4).Range.Text = "23÷2=11, 1"
$t.Cell(17, 5).Range.Text = "82÷6=13, 4"
